$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 9,20

$arr[0,0] = "ECs"
$arr[0,1] = "Efnb1"
$arr[0,2] = "Ephb4"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 10.31211433333333
$arr[0,7] = 30.936343
$arr[0,8] = 0.633340936097251
$arr[0,9] = 0.633340936097251
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 23.746319
$arr[0,13] = 71.238957
$arr[0,14] = 0.7135031414879517
$arr[0,15] = 0.7135031414879517
$arr[0,16] = 244.8747565238057
$arr[0,17] = 2203.872808714251
$arr[0,18] = 0.4518907475383086
$arr[0,19] = 0.4518907475383086

$arr[1,0] = "ECs"
$arr[1,1] = "Efnb1"
$arr[1,2] = "Ephb4"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 10.31211433333333
$arr[1,7] = 30.936343
$arr[1,8] = 0.633340936097251
$arr[1,9] = 0.633340936097251
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 4.865208333333334
$arr[1,13] = 14.595625
$arr[1,14] = 0.1461844014571983
$arr[1,15] = 0.1461844014571983
$arr[1,16] = 50.17058458881945
$arr[1,17] = 451.5352612993751
$arr[1,18] = 0.09258456566171831
$arr[1,19] = 0.09258456566171828

$arr[2,0] = "ECs"
$arr[2,1] = "Efnb1"
$arr[2,2] = "Ephb4"
$arr[2,3] = "sCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 10.31211433333333
$arr[2,7] = 30.936343
$arr[2,8] = 0.633340936097251
$arr[2,9] = 0.633340936097251
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 4.669782333333333
$arr[2,13] = 14.009347
$arr[2,14] = 0.1403124570548501
$arr[2,15] = 0.1403124570548501
$arr[2,16] = 48.15532933311344
$arr[2,17] = 433.397963998021
$arr[2,18] = 0.08886562289722406
$arr[2,19] = 0.08886562289722406

$arr[3,0] = "FAPs"
$arr[3,1] = "Efnb1"
$arr[3,2] = "Ephb4"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 4.103438
$arr[3,7] = 12.310314
$arr[3,8] = 0.2520215719230645
$arr[3,9] = 0.2520215719230645
$arr[3,10] = 2
$arr[3,11] = 0.6666666666666666
$arr[3,12] = 23.746319
$arr[3,13] = 71.238957
$arr[3,14] = 0.7135031414879517
$arr[3,15] = 0.7135031414879517
$arr[3,16] = 97.44154774472199
$arr[3,17] = 876.9739297024979
$arr[3,18] = 0.1798181832898383
$arr[3,19] = 0.1798181832898383

$arr[4,0] = "FAPs"
$arr[4,1] = "Efnb1"
$arr[4,2] = "Ephb4"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 4.103438
$arr[4,7] = 12.310314
$arr[4,8] = 0.2520215719230645
$arr[4,9] = 0.2520215719230645
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 4.865208333333334
$arr[4,13] = 14.595625
$arr[4,14] = 0.1461844014571983
$arr[4,15] = 0.1461844014571983
$arr[4,16] = 19.96408075291667
$arr[4,17] = 179.67672677625
$arr[4,18] = 0.03684162264587544
$arr[4,19] = 0.03684162264587542

$arr[5,0] = "FAPs"
$arr[5,1] = "Efnb1"
$arr[5,2] = "Ephb4"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 4.103438
$arr[5,7] = 12.310314
$arr[5,8] = 0.2520215719230645
$arr[5,9] = 0.2520215719230645
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 4.669782333333333
$arr[5,13] = 14.009347
$arr[5,14] = 0.1403124570548501
$arr[5,15] = 0.1403124570548501
$arr[5,16] = 19.16216227832867
$arr[5,17] = 172.459460504958
$arr[5,18] = 0.03536176598735079
$arr[5,19] = 0.03536176598735079

$arr[6,0] = "sCs"
$arr[6,1] = "Efnb1"
$arr[6,2] = "Ephb4"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1.866538
$arr[6,7] = 5.599614
$arr[6,8] = 0.1146374919796846
$arr[6,9] = 0.1146374919796846
$arr[6,10] = 2
$arr[6,11] = 0.6666666666666666
$arr[6,12] = 23.746319
$arr[6,13] = 71.238957
$arr[6,14] = 0.7135031414879517
$arr[6,15] = 0.7135031414879517
$arr[6,16] = 44.323406773622
$arr[6,17] = 398.910660962598
$arr[6,18] = 0.08179421065980483
$arr[6,19] = 0.08179421065980481

$arr[7,0] = "sCs"
$arr[7,1] = "Efnb1"
$arr[7,2] = "Ephb4"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1.866538
$arr[7,7] = 5.599614
$arr[7,8] = 0.1146374919796846
$arr[7,9] = 0.1146374919796846
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 4.865208333333334
$arr[7,13] = 14.595625
$arr[7,14] = 0.1461844014571983
$arr[7,15] = 0.1461844014571983
$arr[7,16] = 9.081096232083334
$arr[7,17] = 81.72986608875001
$arr[7,18] = 0.01675821314960456
$arr[7,19] = 0.01675821314960456

$arr[8,0] = "sCs"
$arr[8,1] = "Efnb1"
$arr[8,2] = "Ephb4"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 1.866538
$arr[8,7] = 5.599614
$arr[8,8] = 0.1146374919796846
$arr[8,9] = 0.1146374919796846
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 4.669782333333333
$arr[8,13] = 14.009347
$arr[8,14] = 0.1403124570548501
$arr[8,15] = 0.1403124570548501
$arr[8,16] = 8.716326176895333
$arr[8,17] = 78.446935592058
$arr[8,18] = 0.01608506817027522
$arr[8,19] = 0.01608506817027521

$ws.Range("A2:T10").Value = $arr
